# Apply the stock-report adjustments described by the commit diff.
# Quantities (col F) and their derived stock values (col G) are corrected
# downward/adjusted for a number of SKUs, their Sub Total rows (col B) are
# recomputed, and a handful of duplicate-SKU row pairs have their
# Code/MRP/Qty/Value (B/E/F/G) swapped between the two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 53
$ws.Range("G6").Value = 1583.64
$ws.Range("F7").Value = 93
$ws.Range("G7").Value = 4355.19
$ws.Range("B10").Value = 26374.44
$ws.Range("F22").Value = 57
$ws.Range("G22").Value = 422.37
$ws.Range("F31").Value = 28
$ws.Range("G31").Value = 988.96
$ws.Range("B32").Value = 12225.35
$ws.Range("F58").Value = 20
$ws.Range("G58").Value = 3440.8
$ws.Range("F83").Value = 100
$ws.Range("G83").Value = 15067
$ws.Range("F84").Value = 20
$ws.Range("G84").Value = 2049.2
$ws.Range("B90").Value = 166677.49
$ws.Range("F115").Value = 181
$ws.Range("G115").Value = 17522.61
$ws.Range("B117").Value = 11363.67
$ws.Range("F136").Value = 16
$ws.Range("G136").Value = 1200.32
$ws.Range("B138").Value = 2257.71
$ws.Range("F144").Value = 950
$ws.Range("G144").Value = 8027.5
$ws.Range("F146").Value = 16
$ws.Range("G146").Value = 1347.04
$ws.Range("B147").Value = 12458.68
$ws.Range("F149").Value = 213
$ws.Range("G149").Value = 13802.4
$ws.Range("F150").Value = 24
$ws.Range("G150").Value = 1115.76
$ws.Range("F152").Value = 57
$ws.Range("G152").Value = 5032.53
$ws.Range("B156").Value = 28822.88
$ws.Range("F187").Value = 19
$ws.Range("G187").Value = 949.4299999999999
$ws.Range("F190").Value = 1
$ws.Range("G190").Value = 82.01000000000001
$ws.Range("F197").Value = 16
$ws.Range("G197").Value = 992.96
$ws.Range("F211").Value = 41
$ws.Range("G211").Value = 4149.2
$ws.Range("F212").Value = 4
$ws.Range("G212").Value = 1868.76
$ws.Range("F214").Value = 37
$ws.Range("G214").Value = 3244.9
$ws.Range("B216").Value = 32668.13
$ws.Range("B219").Value = 63565
$ws.Range("E219").Value = 109.19
$ws.Range("F219").Value = 60
$ws.Range("G219").Value = 6162.6
$ws.Range("B220").Value = 61610
$ws.Range("E220").Value = 122.71
$ws.Range("F220").Value = -58
$ws.Range("G220").Value = -5957.18
$ws.Range("F225").Value = 72
$ws.Range("G225").Value = 8224.559999999999
$ws.Range("F231").Value = 2
$ws.Range("G231").Value = 236.92
$ws.Range("B232").Value = 55356
$ws.Range("E232").Value = 54.04
$ws.Range("F232").Value = -158
$ws.Range("G232").Value = -7527.12
$ws.Range("B233").Value = 63510
$ws.Range("E233").Value = 50.66
$ws.Range("F233").Value = 113
$ws.Range("G233").Value = 5383.32
$ws.Range("F237").Value = 5
$ws.Range("G237").Value = 1513.95
$ws.Range("F250").Value = 4
$ws.Range("G250").Value = 1977.52
$ws.Range("F255").Value = 510
$ws.Range("G255").Value = 87378.3
$ws.Range("F256").Value = 257
$ws.Range("G256").Value = 38850.69
$ws.Range("B260").Value = 166525.81
$ws.Range("F278").Value = 7
$ws.Range("G278").Value = 961.24
$ws.Range("F280").Value = 129
$ws.Range("G280").Value = 21819.06
$ws.Range("F288").Value = 34
$ws.Range("G288").Value = 3161.66
$ws.Range("F291").Value = 105
$ws.Range("G291").Value = 4516.05
$ws.Range("F292").Value = 38
$ws.Range("G292").Value = 3164.26
$ws.Range("F293").Value = 26
$ws.Range("G293").Value = 1828.32
$ws.Range("F296").Value = 33
$ws.Range("G296").Value = 699.6
$ws.Range("F302").Value = 32
$ws.Range("G302").Value = 6748.48
$ws.Range("B304").Value = 163251.18
$ws.Range("F320").Value = 37
$ws.Range("G320").Value = 2540.05
$ws.Range("F326").Value = 60
$ws.Range("G326").Value = 1784.4
$ws.Range("F328").Value = 32
$ws.Range("G328").Value = 1190.72
$ws.Range("B330").Value = 25286.42
$ws.Range("F334").Value = 188
$ws.Range("G334").Value = 9742.16
$ws.Range("F336").Value = 19
$ws.Range("G336").Value = 829.35
$ws.Range("B346").Value = 23597.09
$ws.Range("F357").Value = 5
$ws.Range("G357").Value = 1306.5
$ws.Range("B358").Value = 34516.18
$ws.Range("B375").Value = 64927
$ws.Range("E375").Value = 17.26
$ws.Range("F375").Value = 106
$ws.Range("G375").Value = 1719.32
$ws.Range("B376").Value = 45718
$ws.Range("E376").Value = 19.38
$ws.Range("F376").Value = -294
$ws.Range("G376").Value = -4768.68
$ws.Range("B385").Value = 53595
$ws.Range("E385").Value = 17.61
$ws.Range("F385").Value = -335
$ws.Range("G385").Value = -4934.55
$ws.Range("B386").Value = 65067
$ws.Range("E386").Value = 15.65
$ws.Range("F386").Value = 126
$ws.Range("G386").Value = 1855.98
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("F454").Value = 46
$ws.Range("G454").Value = 1570.9
$ws.Range("B460").Value = 12109.12
$ws.Range("B463").Value = 60025
$ws.Range("E463").Value = 37.22
$ws.Range("F463").Value = -98
$ws.Range("G463").Value = -3217.34
$ws.Range("B464").Value = 64833
$ws.Range("E464").Value = 34.9
$ws.Range("F464").Value = 95
$ws.Range("G464").Value = 3118.85
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79
$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 105
$ws.Range("G474").Value = 3447.15
$ws.Range("F477").Value = 5
$ws.Range("G477").Value = 226.7
$ws.Range("B478").Value = 226.7
$ws.Range("F509").Value = 191
$ws.Range("G509").Value = 15352.58
$ws.Range("B510").Value = 20757.46
$ws.Range("F555").Value = 14
$ws.Range("G555").Value = 973.84
$ws.Range("B560").Value = 3477.4
$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 18
$ws.Range("G572").Value = 735.66
$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 6
$ws.Range("G573").Value = 245.22
$ws.Range("F578").Value = 55
$ws.Range("G578").Value = 2743.95
$ws.Range("F580").Value = 45
$ws.Range("G580").Value = 2564.55
$ws.Range("F582").Value = 22
$ws.Range("G582").Value = 1253.78
$ws.Range("B583").Value = 12504.42
$ws.Range("F599").Value = 1314
$ws.Range("G599").Value = 214326.54
$ws.Range("F601").Value = 358
$ws.Range("G601").Value = 101267.46
$ws.Range("F602").Value = 313
$ws.Range("G602").Value = 45275.45
$ws.Range("B606").Value = 361717.5
$ws.Range("F613").Value = 126
$ws.Range("G613").Value = 20054.16
$ws.Range("B618").Value = 41312.67
$ws.Range("B619").Value = 1566381.86
$ws.Range("B620").Value = 1566381.86
